$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new test case row was added to the bottom of the script master table
# (row 37): MILESTONE=15, TESTCASE="testT4149", SCRIPT_ITERATION=1,
# EXECUTE_FLAG="YES" — inserted right after the existing testT4245 row.
$ws.Range("A37").Value = 15
$ws.Range("B37").Value = "testT4149"
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = "YES"

# Match the formatting already used throughout the table: column A
# (MILESTONE) is centered both horizontally and vertically, column C
# (SCRIPT_ITERATION) is centered horizontally.
$ws.Range("A37").HorizontalAlignment = -4108
$ws.Range("A37").VerticalAlignment = -4108
$ws.Range("C37").HorizontalAlignment = -4108

# Scroll the view down a bit and move the active selection, matching the
# author's final cursor position after adding the new row.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A27").Select() | Out-Null
